$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# B2 = ja_JP column: replace the untranslated placeholder text with the
# actual Japanese translation.
$ws.Range("B2").Value = "ラヴァは、ハイビスカスへの態度を改めようとするが失敗に終わる。その後、クオーラに刺激されて姉に謝ろうと決心したが、それより先にハイビスカスから詫びの手紙があった。`n"

# C2 = en_US column: replace the untranslated placeholder text with the
# actual English translation.
$ws.Range("C2").Value = "Lava attempts to manage her attitude towards Hibiscus, but fails in the end. Provoked by Cuora, she looks to apologize to her sister, but is beaten to the punch by a letter from Hibiscus.`n"

# The multi-line translations trigger an automatic row-height adjustment;
# restore row 2 to its original (un-customized) auto-fit height so the
# serialized row element doesn't pick up a stray ht/customHeight attribute.
$ws.Rows.Item(2).AutoFit()
